# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K") values for rows 2-59 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 5
    3  = 1
    4  = 1
    5  = 2
    6  = 2
    7  = 2
    8  = 0
    9  = 0
    10 = 0
    11 = 2
    12 = 0
    13 = 0
    14 = 2
    15 = 1
    16 = 2
    17 = 2
    18 = 2
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 0
    25 = 2
    26 = 0
    27 = 1
    28 = 1
    29 = 3
    30 = 0
    31 = 2
    32 = 1
    33 = 2
    34 = 1
    35 = 0
    36 = 0
    37 = 2
    38 = 0
    39 = 2
    40 = 1
    41 = 1
    42 = 0
    43 = 1
    44 = 1
    45 = 1
    46 = 1
    47 = 3
    48 = 1
    49 = 0
    50 = 2
    51 = 3
    52 = 1
    53 = 0
    54 = 0
    55 = 1
    56 = 1
    57 = 1
    58 = 0
    59 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
